$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above the old row 10 ("curriculo"/...), pushing
# everything from the old row 10 downward by two rows.
$ws.Rows("10:11").Insert()

# The newly inserted rows come back blank with a generic style; clone the
# banded formatting (fill/border/number format) from row 9 - the last row
# of the style "group" that rows 10-11 should now belong to - so the new
# rows 10-11 keep the same style indices as rows 2-9 (fill 4 / border 1).
$ws.Range("A9:J9").Copy()
$ws.Range("A10:J11").PasteSpecial(-4122)

# Row 10: new project "fakesensor"; Row 11: new project "pontilhar"
# (project names are written before descriptions so new shared strings
# land in the same order the source workbook has them in)
$ws.Range("A10").Value = "fakesensor"
$ws.Range("A11").Value = "pontilhar"
$ws.Range("B10").Value = "FKS é uma boa maneira de simular um sensor que produz arquivos de dados"
$ws.Range("B11").Value = "Aplicativo web para criação de desenhos com pontos coloridos"

$ws.Range("C10").Value = 42928
$ws.Range("D10").Value = "DEV"
$ws.Range("E10").Value = "X"
$ws.Range("I10").Value = "MPS"

$ws.Range("C11").Value = 42929
$ws.Range("D11").Value = "DEV"
$ws.Range("E11").Value = "X"
$ws.Range("I11").Value = "MPS"

# Row 9 ("elk") now also has an "Externos" (X) mark.
$ws.Range("J9").Value = "X"

# Update selection: Excel no longer needs to scroll (no topLeftCell),
# the active cell moves to C11.
$ws.Range("C11").Select()
